$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VEGFA165_NRP1")

# --- Add Shobhan, 2023 data point pair (F2/G2 label+value, G3 second value) ---
$ws.Range("F2").Value = "Shobhan, 2023"
$ws.Range("G2").Value = 5.29
$ws.Range("G3").Value = 7.43

# --- Insert a new row at 7 (pushes the existing "Gu et al., 2002" row down to row 8) ---
$ws.Rows.Item(7).Insert(-4121)

# Give the new row 7 the same look as row 6 (continuing interior-row borders)
$ws.Range("A6:D6").Copy()
$ws.Range("A7:D7").PasteSpecial(-4122)

# Restore the "Gu et al., 2002" values into row 7 (same reference, now mid-table)
$ws.Range("A7").Value = "Gu et al., 2002"
$ws.Range("B7").Value = "Radioligand"
$ws.Range("C7").Value = 0.93
$ws.Range("D7").Value = 0.71

# Row 8 keeps the old "Gu et al., 2002" formatting (thick bottom border/closing row);
# turn it into the new "Unpublished data" / Shobhan average+SE summary row.
$ws.Range("A8").Value = "Unpublished data"
$ws.Range("B8").Value = "SPR"
$ws.Range("C8").Formula = "=AVERAGE(G2:G3)"
$ws.Range("D8").Formula = "=STDEVA(G2:G3)/SQRT(2)"

$ws.Range("A13").Select()
